$d = $word.ActiveDocument

# Remember where the existing content ends.
$startCount = $d.Paragraphs.Count

# --- Step 1: create all six new paragraphs (as plain, unformatted
# paragraphs) first. Doing the paragraph-break insertion before any bold
# text is typed keeps later paragraphs from inheriting bold formatting
# from the heading line. ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # startCount+1 : blank
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # startCount+2 : "October 25th:"
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # startCount+3 : blank (bold)
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # startCount+4 : "Outlined architectural design..."
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # startCount+5 : blank
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # startCount+6 : "Began working on features in react. "

$pBlank1    = $d.Paragraphs.Item($startCount + 1)
$pHeading   = $d.Paragraphs.Item($startCount + 2)
$pBlankBold = $d.Paragraphs.Item($startCount + 3)
$pOutline   = $d.Paragraphs.Item($startCount + 4)
$pBlank2    = $d.Paragraphs.Item($startCount + 5)
$pReact     = $d.Paragraphs.Item($startCount + 6)

# Helper trick: a paragraph that should stay empty still needs its
# w:rPr (the paragraph-mark run properties) set. Typing a placeholder
# character, applying formatting, then deleting the character again
# leaves the paragraph mark's properties set without leaving a stray
# empty <w:r> behind.
function Set-BlankParaFormat($para, [bool]$bold) {
    $r = $para.Range
    $r.InsertAfter("X")
    $r.Font.Name = "Times New Roman"
    $r.Font.NameBi = "Times New Roman"
    if ($bold) {
        $r.Font.Bold = 1
        $r.Font.BoldBi = 1
    }
    $delR = $d.Range($para.Range.Start, $para.Range.Start + 1)
    $delR.Delete()
}

# --- blank paragraph right after "Worked on scenarios..." ---
Set-BlankParaFormat $pBlank1 $false

# --- "October 25th:" heading (bold, centered, three runs) ---
$r = $pHeading.Range
$r.InsertAfter("October 25")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"
$r.Font.Bold = 1
$r.Font.BoldBi = 1

$r = $d.Range($pHeading.Range.End - 1, $pHeading.Range.End - 1)
$r.InsertAfter("th")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"
$r.Font.Bold = 1
$r.Font.BoldBi = 1
$r.Font.Superscript = $true

$r = $d.Range($pHeading.Range.End - 1, $pHeading.Range.End - 1)
$r.InsertAfter(":")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"
$r.Font.Bold = 1
$r.Font.BoldBi = 1

# --- blank paragraph right after the heading keeps the bold mark ---
Set-BlankParaFormat $pBlankBold $true

# --- "Outlined architectural design, prototype components, and
# selected features." (four separate runs, not bold) ---
$r = $pOutline.Range
$r.InsertAfter("Outline")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"

$r = $d.Range($pOutline.Range.End - 1, $pOutline.Range.End - 1)
$r.InsertAfter("d")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"

$r = $d.Range($pOutline.Range.End - 1, $pOutline.Range.End - 1)
$r.InsertAfter(" architectural design")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"

$r = $d.Range($pOutline.Range.End - 1, $pOutline.Range.End - 1)
$r.InsertAfter(", prototype components, and selected features.")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"

# --- blank paragraph stays plain ---
Set-BlankParaFormat $pBlank2 $false

# --- "Began working on features in react. " ---
$r = $pReact.Range
$r.InsertAfter("Began working on features in react. ")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"
